$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.702409863471985
$ws.Range("B1").Value = 2.216505289077759
$ws.Range("C1").Value = 2.246217250823975
$ws.Range("D1").Value = 7.308047771453857
$ws.Range("E1").Value = 0.7607421278953552
